# Project Sample Project is saved.TEST Author: admin. Type: SAVE.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell B11 held the text "R40" (the row-40 label); update it to the literal
# text "1" while leaving the cell's existing style/number format untouched.
#
# A plain `Range.Value = "1"` assignment would be auto-coerced to the number
# 1 (it "looks like" a number), which is not what we want here — the target
# cell must keep storing a text value. Stage the text in a scratch cell
# that's explicitly formatted as Text ("@") so it is not reinterpreted as a
# number, then copy/paste *values only* into B11 so its own formatting
# (style) is left exactly as it was.
$scratch = $ws.Range("Z1")
$scratch.NumberFormat = "@"
$scratch.Value = "1"
$scratch.Copy()
$ws.Range("B11").PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false
$scratch.Clear()
